# repull data, push all data, mean calculation
# Update column F (dSF) values on Sheet1 for the rows that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    5  = -5
    6  = 0
    10 = -7
    12 = 2
    13 = -1
    17 = -5
    18 = 0
    19 = -3
    22 = 3
    25 = 1
    28 = -3
    30 = -4
    34 = 6
    35 = -4
    39 = -2
    43 = 4
    45 = -7
    48 = -1
    49 = -4
    50 = -10
    52 = -6
    55 = -6
    57 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
